# Backup QR Scanner data - 4/5/2025, 9:52:50 PM
$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet so it lands at the end
# of the tab order (Worksheets.Add() with no args inserts before the active
# sheet, which is not what we want here).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Mazinjsbdbmsndbd"

# Header row
$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Student ID"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Log Date"
$ws.Range("E1").Value = "Log Time"

# Force text storage for columns that look numeric/date-like (Student ID,
# Log Date, Log Time), matching how the QR-scanner logger writes every
# other sheet in this workbook -- these are plain text, not real
# Excel numbers/dates.
$textCols = "B", "C", "D", "E"
foreach ($col in $textCols) {
    $ws.Range($col + "1:" + $col + "5").NumberFormat = "@"
}

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "231249"
$ws.Range("C2").Value = "Mazinjsbdbmsndbd"
$ws.Range("D2").Value = "2025-04-05"
$ws.Range("E2").Value = "21:52:35"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "231249"
$ws.Range("C3").Value = "Mazinjsbdbmsndbd"
$ws.Range("D3").Value = "2025-04-05"
$ws.Range("E3").Value = "21:52:39"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "231249"
$ws.Range("C4").Value = "Mazinjsbdbmsndbd"
$ws.Range("D4").Value = "2025-04-05"
$ws.Range("E4").Value = "21:52:42"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "231249"
$ws.Range("C5").Value = "Mazinjsbdbmsndbd"
$ws.Range("D5").Value = "2025-04-05"
$ws.Range("E5").Value = "21:52:45"
